# Add a new orchard record ("Beilke Family Farm") as row 48 of Sheet1.
# Columns: A=id, B=orchardName, C=Latitude, D=Longitude, E=City, F=State,
#          G=Breeds, H=Address, I=GMapsLink, J=PhoneNumber, K=Website, L=Email

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(48, 2).Value = "Beilke Family Farm"
$ws.Cells.Item(48, 3).Value = 45.050645267309697
$ws.Cells.Item(48, 4).Value = -122.954332824716
$ws.Cells.Item(48, 5).Value = "Salem"
$ws.Cells.Item(48, 6).Value = "Oregon"
$ws.Cells.Item(48, 7).Value = "['Gala', 'Honeycrisp', 'Ruby Mac', 'Smoothie', 'Jonagold']"
$ws.Cells.Item(48, 8).Value = "4925 Rockdale St NE, Salem, OR 97305"
$ws.Cells.Item(48, 9).Value = "https://www.google.com/maps/place/Beilke+Family+Farm/@45.0487481,-122.9589111,13.96z/data=!4m6!3m5!1s0x549557975f5f8607:0xd2ca29d38342900!8m2!3d45.049593!4d-122.9571511!16s%2Fg%2F11b5pjd1k5!5m1!1e3?entry=ttu&g_ep=EgoyMDI0MDkwMi4wIKXMDSoASAFQAw%3D%3D"
$ws.Cells.Item(48, 10).Value = "(503) 393 1077"
$ws.Cells.Item(48, 11).Value = "https://www.beilkefamilyfarm.com/"
$ws.Cells.Item(48, 12).Value = "beilkeff@gmail.com"

# Match the workbook's existing style for the GMapsLink / Website / Email
# columns (fill alignment, same as other rows).
$ws.Cells.Item(48, 9).HorizontalAlignment = 5
$ws.Cells.Item(48, 11).HorizontalAlignment = 5
$ws.Cells.Item(48, 12).HorizontalAlignment = 5

# Reflect the author's final scroll position / active cell selection.
$ws.Activate()
$ws.Range("B49").Select()
